{"js": "// The author split the sentence \"\uc704\uc640 \uac19\uc740 \uc774\uc720\ub85c \uc81c\uc678.\" into three runs,\n// effectively rewording it to \"\uc704\uc640 \uac19\uc740 \uc774\uc720\uc640 \uac01 \ub098\ub77c\uc758 \uaddc\uaca9\uc774 \uc11c\ub85c\n// \uc720\uc0ac\ud558\uae30 \ub54c\ubb38\uc5d0 \uc81c\uc678\ud568.\" Locate the run by its original text and\n// replace its contents with the new wording (formatting is preserved\n// because insertText(\"Replace\") rewrites the text in place on the\n// matched range).\nconst body = context.document.body;\n\nconst results = body.search(\"\uc704\uc640 \uac19\uc740 \uc774\uc720\ub85c \uc81c\uc678.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \uc704\uc640 \uac19\uc740 \uc774\uc720\ub85c \uc81c\uc678.\");\n}\n\nresults.items[0].insertText(\n  \"\uc704\uc640 \uac19\uc740 \uc774\uc720\uc640 \uac01 \ub098\ub77c\uc758 \uaddc\uaca9\uc774 \uc11c\ub85c \uc720\uc0ac\ud558\uae30 \ub54c\ubb38\uc5d0 \uc81c\uc678\ud568.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# The author reworded the bullet \"\uc81c\ud488 \uaddc\uaca9 \uae30\uc900\uad6d -> \uc704\uc640 \uac19\uc740 \uc774\uc720\ub85c\n# \uc81c\uc678.\" so that the trailing sentence reads \"\uc704\uc640 \uac19\uc740 \uc774\uc720\uc640 \uac01 \ub098\ub77c\uc758\n# \uaddc\uaca9\uc774 \uc11c\ub85c \uc720\uc0ac\ud558\uae30 \ub54c\ubb38\uc5d0 \uc81c\uc678\ud568.\" instead of \"\uc704\uc640 \uac19\uc740 \uc774\uc720\ub85c \uc81c\uc678.\"\n# Find the exact original sentence and replace it with the new wording,\n# leaving the \"\uc81c\ud488 \uaddc\uaca9 \uae30\uc900\uad6d -> \" lead-in and paragraph/list formatting\n# untouched.\n\n$d = $word.ActiveDocument\n\n$oldText = \"\uc704\uc640 \uac19\uc740 \uc774\uc720\ub85c \uc81c\uc678.\"\n$newText = \"\uc704\uc640 \uac19\uc740 \uc774\uc720\uc640 \uac01 \ub098\ub77c\uc758 \uaddc\uaca9\uc774 \uc11c\ub85c \uc720\uc0ac\ud558\uae30 \ub54c\ubb38\uc5d0 \uc81c\uc678\ud568.\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\n    $oldText,   # FindText\n    $true,      # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap: wdFindContinue\n    $false,     # Format\n    $newText,   # ReplaceWith\n    1           # Replace: wdReplaceOne\n)\n\nif (-not $found) {\n    throw \"Target sentence not found: $oldText\"\n}\n\n# Defensively keep the run's original Korean font in place (it already\n# carries over from the matched run, this just makes it explicit).\n$after = $d.Content\n$after.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0) | Out-Null\n$after.Font.NameAscii = \"NanumGothic\"\n$after.Font.NameFarEast = \"NanumGothic\"\n$after.Font.Name = \"NanumGothic\"\n"}
